$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

# Update achievement cost requirements
# Big Spender (row 23): 10 -> 5
$ws.Range("D23").Value = 5
# Endless Wallet (row 24): 25 -> 10
$ws.Range("D24").Value = 10
# Fat Cat (row 25): 50 -> 20
$ws.Range("D25").Value = 20

# Update the active selection to D25
$ws.Range("D25").Select()
